$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("A7").Value = "Максим Шило"
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = "mamkotraxer@gmail.com"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "380984514236"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "Mamu ebal"
